$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "[Python] mAP(mean Average Precision) 예시 및 코드"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-mAPmean-Average-Precision-%EC%98%88%EC%8B%9C-%EB%B0%8F-%EC%BD%94%EB%93%9C"

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D37").Value = "[Paper Review] CaSS : A Channel-aware Self-supervised Representation Learning for Multivariate Time Series Classification"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=2000&mod=document&pageid=1"

$ws.Range("D46").Value = "[Bioinformatics] 2022년 07월, 제16회 통계유전학워크샵 [한국유전체학회]"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/477"

$ws.Range("D51").Value = "[pyside6] 창 타이틀바에서 닫기 버튼, 최소화 버튼, 최대화 버튼 등을 비활성화 되게 하려면?"
$ws.Range("E51").Value = "https://bskyvision.com/1288"
